# Updated cryptos list on Mon Apr 10 12:17:06 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as "28.544.38" or "0.00001101" that Excel would
# otherwise auto-convert to a number (dropping trailing zeros / using scientific
# notation). Force those specific cells to Text format first so the literal string
# is preserved exactly, matching the source data. Cells that do not change are left
# completely untouched (no format applied).
$dRanges = @("D2:D9", "D11:D13", "D15:D31", "D33:D35", "D37", "D39:D42", "D44:D48", "D50:D51")
foreach ($r in $dRanges) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.544.38"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.873.82"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -1.93%  "
$ws.Range("D5").Value = "315.96"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("D7").Value = "0.5095"
$ws.Range("E7").Value = "  -1.70%  "
$ws.Range("D8").Value = "0.3904"
$ws.Range("E8").Value = "  -2.23%  "
$ws.Range("D9").Value = "0.08366"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").Value = "41.92"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "6.215"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "1.873.56"
$ws.Range("E13").Value = "  +4.71%  "
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "7.271"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "91.29"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "0.06758"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "17.73"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "5.926"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "28.589.81"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").Value = "2.211"
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("D26").Value = "2.083.67"
$ws.Range("E26").Value = "  +3.95%  "
$ws.Range("D27").Value = "157.68"
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").Value = "20.61"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").Value = "2.417"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "126.23"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "0.1039"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "5.736"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "3.621"
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "0.02461"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "8.926"
$ws.Range("E37").Value = "  -3.54%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").Value = "5.057"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "1.182"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").Value = "1.238"
$ws.Range("E41").Value = "  -2.72%  "
$ws.Range("D42").Value = "0.6372"
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").Value = "1.009"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "0.6007"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").Value = "13.08"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "3.689"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "2.006"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").Value = "122.64"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "0.06815"
$ws.Range("E51").Value = "  -1.21%  "
